$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7144781
$ws.Range("J19").Value = 10002018
$ws.Range("L19").Value = 10002018
$ws.Range("N19").Value = -10002368
$ws.Range("H33").Value = 269.41177
$ws.Range("I33").Value = 258.1111
$ws.Range("K33").Value = 258.1111
$ws.Range("M33").Value = -29.11110000000002
$ws.Range("H53").Value = 698.7727
$ws.Range("I53").Value = 602.25
$ws.Range("J53").Value = 753.9286
$ws.Range("K53").Value = 602.25
$ws.Range("L53").Value = 753.9286
$ws.Range("M53").Value = 34.75
$ws.Range("N53").Value = -2027.9286
$ws.Range("H64").Value = 3950.4
$ws.Range("J64").Value = 3936.5
$ws.Range("L64").Value = 3936.5
$ws.Range("N64").Value = -4432.5
$ws.Range("H67").Value = 3950.4
$ws.Range("J67").Value = 3936.5
$ws.Range("L67").Value = 3936.5
$ws.Range("N67").Value = -5652.5
$ws.Range("H106").Value = 8125.2856
$ws.Range("I106").Value = 8027.636
$ws.Range("J106").Value = 8483.333000000001
$ws.Range("K106").Value = 8027.636
$ws.Range("L106").Value = 8483.333000000001
$ws.Range("M106").Value = -7396.636
$ws.Range("N106").Value = -9745.333000000001
$ws.Range("H107").Value = 906.85187
$ws.Range("J107").Value = 1821
$ws.Range("L107").Value = 1821
$ws.Range("N107").Value = -5661
$ws.Range("H137").Value = 1835.875
$ws.Range("I137").Value = 1843.5
$ws.Range("K137").Value = 5530.5
$ws.Range("M137").Value = -2980.5
$ws.Range("H138").Value = 6217.2354
$ws.Range("J138").Value = 7013.4346
$ws.Range("L138").Value = 21040.3038
$ws.Range("N138").Value = -31320.3038

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2872.2
$ws.Range("I45").Value = 1590.25
$ws.Range("K45").Value = 1590.25
$ws.Range("M45").Value = -1213.25
$ws.Range("H46").Value = 18000.334
$ws.Range("I46").Value = 12999
$ws.Range("J46").Value = 18455
$ws.Range("K46").Value = 12999
$ws.Range("L46").Value = 18455
$ws.Range("M46").Value = -12680
$ws.Range("N46").Value = -19093

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 34993
$ws.Range("J21").Value = 34993
$ws.Range("L21").Value = 34993
$ws.Range("N21").Value = -35465
$ws.Range("I22").Value = 697
$ws.Range("J22").Value = 2160
$ws.Range("K22").Value = 697
$ws.Range("L22").Value = 2160
$ws.Range("M22").Value = -524
$ws.Range("N22").Value = -2506
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H107").Value = 4769
$ws.Range("I107").Value = 4879.625
$ws.Range("K107").Value = 4879.625
$ws.Range("M107").Value = -2959.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23812438
$ws.Range("I31").Value = 34485420
$ws.Range("J31").Value = 3482.5386
$ws.Range("K31").Value = 34485420
$ws.Range("L31").Value = 3482.5386
$ws.Range("M31").Value = -34485125
$ws.Range("N31").Value = -4072.5386
$ws.Range("H34").Value = 23812438
$ws.Range("I34").Value = 34485420
$ws.Range("J34").Value = 3482.5386
$ws.Range("K34").Value = 34485420
$ws.Range("L34").Value = 3482.5386
$ws.Range("M34").Value = -34485218
$ws.Range("N34").Value = -3886.5386
$ws.Range("H63").Value = 79271
$ws.Range("J63").Value = 79271
$ws.Range("L63").Value = 79271
$ws.Range("N63").Value = -80643
$ws.Range("H66").Value = 79271
$ws.Range("J66").Value = 79271
$ws.Range("L66").Value = 237813
$ws.Range("N66").Value = -244677
$ws.Range("H94").Value = 995.2222
$ws.Range("I94").Value = 885.6667
$ws.Range("J94").Value = 1104.7778
$ws.Range("K94").Value = 885.6667
$ws.Range("L94").Value = 1104.7778
$ws.Range("M94").Value = -434.6667
$ws.Range("N94").Value = -2006.7778
$ws.Range("H105").Value = 2380.3333
$ws.Range("I105").Value = 1741.5834
$ws.Range("K105").Value = 1741.5834
$ws.Range("M105").Value = 5.416600000000017
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 106.57143
$ws.Range("I40").Value = 41
$ws.Range("K40").Value = 164
$ws.Range("M40").Value = -95
$ws.Range("H51").Value = 18779
$ws.Range("I51").Value = 8004
$ws.Range("J51").Value = 24166.5
$ws.Range("K51").Value = 24012
$ws.Range("L51").Value = 72499.5
$ws.Range("M51").Value = -23552
$ws.Range("N51").Value = -73419.5
$ws.Range("H68").Value = 3226.3572
$ws.Range("I68").Value = 2795
$ws.Range("J68").Value = 3657.7144
$ws.Range("K68").Value = 8385
$ws.Range("L68").Value = 10973.1432
$ws.Range("M68").Value = -7574
$ws.Range("N68").Value = -12595.1432
$ws.Range("H71").Value = 3226.3572
$ws.Range("I71").Value = 2795
$ws.Range("J71").Value = 3657.7144
$ws.Range("K71").Value = 25155
$ws.Range("L71").Value = 32919.4296
$ws.Range("M71").Value = -21099
$ws.Range("N71").Value = -41031.4296
$ws.Range("H107").Value = 7006002.5
$ws.Range("I107").Value = 2989.6
$ws.Range("J107").Value = 11382885
$ws.Range("K107").Value = 8968.799999999999
$ws.Range("L107").Value = 34148655
$ws.Range("M107").Value = -7048.799999999999
$ws.Range("N107").Value = -34152495
$ws.Range("H121").Value = 5003995
$ws.Range("I121").Value = 250
$ws.Range("J121").Value = 5559967
$ws.Range("K121").Value = 750
$ws.Range("L121").Value = 16679901
$ws.Range("M121").Value = 560
$ws.Range("N121").Value = -16682521
$ws.Range("H129").Value = 17862580
$ws.Range("I129").Value = 22729846
$ws.Range("K129").Value = 68189538
$ws.Range("M129").Value = -68184538
$ws.Range("H132").Value = 2106
$ws.Range("I132").Value = 2088.4
$ws.Range("J132").Value = 2194
$ws.Range("K132").Value = 18795.6
$ws.Range("L132").Value = 19746
$ws.Range("M132").Value = -16265.6
$ws.Range("N132").Value = -24806

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3402.6316
$ws.Range("I102").Value = 3517.4
$ws.Range("K102").Value = 3517.4
$ws.Range("M102").Value = -1895.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1061.3928
$ws.Range("I55").Value = 568.3570999999999
$ws.Range("J55").Value = 1554.4286
$ws.Range("K55").Value = 568.3570999999999
$ws.Range("L55").Value = 1554.4286
$ws.Range("M55").Value = -395.3570999999999
$ws.Range("N55").Value = -1900.4286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 28154.857
$ws.Range("I55").Value = 22271.25
$ws.Range("J55").Value = 35999.668
$ws.Range("K55").Value = 22271.25
$ws.Range("L55").Value = 35999.668
$ws.Range("M55").Value = -21994.25
$ws.Range("N55").Value = -36553.668
$ws.Range("H107").Value = 3402.639
$ws.Range("I107").Value = 1498
$ws.Range("J107").Value = 4135.1924
$ws.Range("K107").Value = 4494
$ws.Range("L107").Value = 12405.5772
$ws.Range("M107").Value = -2574
$ws.Range("N107").Value = -16245.5772
